# Applies the content updates described by the diff:
# - header date changes from 2025-03-03 Monday -> 2025-03-04 Tuesday
# - each "AxB=C" multiplication answer cell is replaced with a new equation

$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-03 Monday", "2025-03-04 Tuesday"),
    @("605×4=2420", "943×9=8487"),
    @("815×2=1630", "828×8=6624"),
    @("905×9=8145", "684×4=2736"),
    @("931×8=7448", "195×3=585"),
    @("715×9=6435", "212×4=848"),
    @("257×5=1285", "522×6=3132"),
    @("584×6=3504", "362×6=2172"),
    @("746×9=6714", "535×9=4815"),
    @("188×8=1504", "590×8=4720"),
    @("106×8=848", "527×5=2635"),
    @("514×2=1028", "477×9=4293"),
    @("712×7=4984", "213×9=1917"),
    @("716×6=4296", "237×9=2133"),
    @("264×8=2112", "711×6=4266"),
    @("571×3=1713", "899×3=2697"),
    @("961×7=6727", "227×2=454"),
    @("800×3=2400", "352×4=1408"),
    @("566×2=1132", "321×9=2889"),
    @("372×4=1488", "953×3=2859"),
    @("659×2=1318", "964×9=8676"),
    @("332×7=2324", "240×5=1200"),
    @("858×2=1716", "287×3=861"),
    @("927×5=4635", "843×6=5058"),
    @("956×7=6692", "696×3=2088"),
    @("151×7=1057", "699×7=4893")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

Write-Host "All replacements attempted."
